# Generate Report for Handoff
# Updates the localization status report to reflect that the content is now
# "Ready for handoff" (previously "In Translation") and refreshes the
# handoff timestamps recorded on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Latest handoff timestamps ------------------------------------------
# zh-cn handoff datetime
$wsZhCn.Range("H2").Value = "2016-08-23 20:41:12"

# de-de handoff datetime (also mirrored on the Overview sheet)
$wsDeDe.Range("H2").Value    = "2016-08-23 20:41:17"
$wsOverview.Range("G2").Value = "2016-08-23 20:41:17"

# --- Column widths -------------------------------------------------------
# Widen the "Status" columns so the longer "Ready for handoff" text fits,
# matching the autofit that Excel performs on edit. (16.3333... character
# units is the closest value the ColumnWidth property can resolve to the
# target display width.)
$wsOverview.Range("E1").ColumnWidth = 16.3333333333333
$wsOverview.Range("F1").ColumnWidth = 16.3333333333333
$wsZhCn.Range("C1").ColumnWidth     = 16.3333333333333
$wsDeDe.Range("C1").ColumnWidth     = 16.3333333333333
